$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 422, shifting existing rows 422-528 down to 423-529
$ws.Rows("422:422").Insert()

# Populate the newly inserted row 422 with the new weekly data entry
$ws.Range("A422").Value = 4
$ws.Range("B422").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C422").Value = "Los Lagos"
$ws.Range("D422").Value = 45173
$ws.Range("E422").Value = 10
$ws.Range("F422").Value = 100112045
$ws.Range("G422").Value = "Zapallo"
$ws.Range("H422").Value = "Paine"
$ws.Range("I422").Value = "1a (guarda)"
$ws.Range("J422").Value = 500
$ws.Range("K422").Value = 800
$ws.Range("L422").Value = 800
$ws.Range("M422").Value = 800
$ws.Range("N422").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O422").Value = "Región de O'Higgins"
$ws.Range("P422").Value = 800
$ws.Range("Q422").Value = 1
$ws.Range("R422").Value = "Hortaliza"
